$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "FPÖ - Freedom Party of Austria (Freiheitliche Partei Österreichs, FPÖ)"
$ws.Range("C1").Value = "GA - The Greens-Green Alternative (Die Grünen-Die Grüne Alternative, GA)"
$ws.Range("D1").Value = "LIF - Liberal Forum (Liberales Forum, LIF)"
$ws.Range("E1").Value = "SPÖ - Social Democratic Party of Austria  (Sozialdemokratische Partei Österreichs, SPÖ)"
$ws.Range("F1").Value = "ÖVP - Austrian People's Party  (Österreichische Volkspartei, ÖVP)"
$ws.Range("G1").Value = "BZÖ - Alliance for the Future of Austria  (Bündnis Zukunft Österreich , BZÖ)"
$ws.Range("H1").Value = "FRANK - Team Frank Stronach (Team Frank Stronach, FRANK)"
$ws.Range("I1").Value = "NEOS - New Austria and Liberal Forum (NEOS Das Neue Österreich und Liberales Forum, NEOS)"
$ws.Range("J1").Value = "PILZ - Peter Pilz List (Liste Peter Pilz, PILZ)"
